$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has headers in B1 ("Jun_13") and C1 ("Jun_10"), with
# data columns A (firm name), B ("UN") and C ("UN", except row 23 which
# holds a special note). Two new reporting-period columns are being added
# in front of the existing "Jun_10" column, so:
#   - insert two new columns at C:D -> the old column C (and its data)
#     shifts right to become column E, carrying its values/format with it
#   - stamp the new column headers: B1 = "Jun_17", C1 = "Jun_15",
#     D1 = "Jun_13" (E1 keeps the old "Jun_10" header automatically)
#   - fill the two freshly inserted data columns (C2:D27) with "UN" to
#     match the existing B/E columns' pattern

$ws.Range("C1:D1").EntireColumn.Insert()

$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

$ws.Range("C2:D27").Value = "UN"

# Match column C's existing custom width (8 characters) on the two new
# columns plus the shifted-right original column.
$ws.Columns("C:E").ColumnWidth = 7.166666666666667

$wb.Save()
